# Update the "Förändrad" (C) column date for every data row, and add the
# friendly display-text argument to every HYPERLINK() formula so it shows
# the "Beteckning" (column A) value instead of the raw URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 211
$hyperlinkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($row = 2; $row -le $lastRow; $row++) {

    # C column: Förändrad date serial 45184 -> 45186
    $ws.Cells.Item($row, 3).Value = 45186

    # Column A holds the "Beteckning" text used as the HYPERLINK friendly name
    $beteckning = $ws.Cells.Item($row, 1).Text

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range($col + $row)
        if ($cell.HasFormula) {
            $oldFormula = $cell.Formula
            if ($oldFormula.StartsWith("=HYPERLINK(") -and -not $oldFormula.Contains(",")) {
                $newFormula = $oldFormula.Substring(0, $oldFormula.Length - 1) + ', "' + $beteckning + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
